$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.717.02"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.231.10"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'253.64"
$ws.Range("E5").Value = "  +8.49%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "'71.62"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "'0.565"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'42.44"
$ws.Range("E10").Value = "  +15.01%  "
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  -5.59%  "
$ws.Range("D12").Value = "'58.61"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'6.94"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "2.556.92"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'15.01"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "'0.860"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "2.236.77"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "41.698.33"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "0.0₃0970"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("D21").Value = "'73.11"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "'2.28"
$ws.Range("E23").Value = "  +17.37%  "
$ws.Range("D24").Value = "'234.90"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'3.75"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "'2.53"
$ws.Range("E27").Value = "  +6.67%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "'170.37"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'20.75"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'0.121"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "'5.52"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "'0.0722"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "'26.60"
$ws.Range("E36").Value = "  +17.82%  "
$ws.Range("D37").Value = "'4.68"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'4.13"
$ws.Range("E38").Value = "  +13.70%  "
$ws.Range("D39").Value = "'0.0285"
$ws.Range("E39").Value = "  +6.38%  "
$ws.Range("B40").Value = "MultiversX"
$ws.Range("C40").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D40").Value = "'70.91"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "'6.03"
$ws.Range("D43").Value = "'0.216"
$ws.Range("E43").Value = "  +13.50%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").Value = "'11.91"
$ws.Range("E44").Value = "  +12.22%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'5.07"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").Value = "'4.84"
$ws.Range("E46").Value = "  +8.19%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.84"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "  +6.25%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.19"
$ws.Range("E51").Value = "  +0.55%  "
